$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $row = $tbl.Rows.Item($i)
    if ($row.Cells.Count -ge 3) {
        $cell = $row.Cells.Item(3)
        $cellText = $cell.Range.Text
        $plainText = $cellText.Trim([char]13, [char]7, ' ', "`t", "`r", "`n")
        if ($plainText.Length -eq 0) {
            $para = $cell.Range.Paragraphs.Item(1)
            $para.Style = "Compact"
        }
    }
}
